$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the shift-preferences text for row 29: "SAT4PM-8PM" -> "SAT 4PM-8PM"
$ws.Range("C29").Value = "FRI 8PM-12AM, SAT 12PM-4PM, SAT 4PM-8PM, SAT 8PM-12AM, SUN 12PM-4PM, SUN 4PM-8PM, SUN 8PM-12AM"

# Column C now holds much longer text, so give it a wide, best-fit-like width
# while columns A:B keep their original width.
$ws.Columns("C").ColumnWidth = 185.75

# Move the selection to C11 and scroll the view so column C is visible
$ws.Range("C11").Select()
$excel.ActiveWindow.ScrollColumn = 3
